{"js": "// Update the worksheet date and all multiplication problems in the table\n// to the new values. Every old value occurs exactly once in the document,\n// so a plain find-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"2024-08-26 Monday\", \"2024-08-27 Tuesday\"],\n  [\"20\u00d741=\", \"45\u00d743=\"],\n  [\"46\u00d780=\", \"52\u00d776=\"],\n  [\"25\u00d771=\", \"93\u00d743=\"],\n  [\"48\u00d730=\", \"72\u00d719=\"],\n  [\"58\u00d799=\", \"71\u00d768=\"],\n  [\"71\u00d732=\", \"65\u00d796=\"],\n  [\"30\u00d724=\", \"39\u00d712=\"],\n  [\"94\u00d728=\", \"43\u00d797=\"],\n  [\"94\u00d770=\", \"73\u00d723=\"],\n  [\"95\u00d761=\", \"11\u00d770=\"],\n  [\"79\u00d741=\", \"70\u00d743=\"],\n  [\"45\u00d783=\", \"45\u00d720=\"],\n  [\"83\u00d722=\", \"45\u00d718=\"],\n  [\"32\u00d781=\", \"90\u00d762=\"],\n  [\"75\u00d721=\", \"17\u00d787=\"],\n  [\"52\u00d782=\", \"89\u00d746=\"],\n  [\"94\u00d763=\", \"26\u00d766=\"],\n  [\"81\u00d789=\", \"22\u00d722=\"],\n  [\"16\u00d735=\", \"37\u00d750=\"],\n  [\"75\u00d754=\", \"90\u00d764=\"],\n  [\"28\u00d755=\", \"17\u00d711=\"],\n  [\"13\u00d788=\", \"95\u00d724=\"],\n  [\"16\u00d780=\", \"32\u00d727=\"],\n  [\"79\u00d791=\", \"19\u00d742=\"],\n  [\"61\u00d737=\", \"26\u00d772=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all multiplication problems in the table\n# to the new values. Each value is unique in the document, so a simple\n# Find/Replace (ReplaceAll) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-26 Monday\", \"2024-08-27 Tuesday\"),\n    @(\"20\u00d741=\", \"45\u00d743=\"),\n    @(\"46\u00d780=\", \"52\u00d776=\"),\n    @(\"25\u00d771=\", \"93\u00d743=\"),\n    @(\"48\u00d730=\", \"72\u00d719=\"),\n    @(\"58\u00d799=\", \"71\u00d768=\"),\n    @(\"71\u00d732=\", \"65\u00d796=\"),\n    @(\"30\u00d724=\", \"39\u00d712=\"),\n    @(\"94\u00d728=\", \"43\u00d797=\"),\n    @(\"94\u00d770=\", \"73\u00d723=\"),\n    @(\"95\u00d761=\", \"11\u00d770=\"),\n    @(\"79\u00d741=\", \"70\u00d743=\"),\n    @(\"45\u00d783=\", \"45\u00d720=\"),\n    @(\"83\u00d722=\", \"45\u00d718=\"),\n    @(\"32\u00d781=\", \"90\u00d762=\"),\n    @(\"75\u00d721=\", \"17\u00d787=\"),\n    @(\"52\u00d782=\", \"89\u00d746=\"),\n    @(\"94\u00d763=\", \"26\u00d766=\"),\n    @(\"81\u00d789=\", \"22\u00d722=\"),\n    @(\"16\u00d735=\", \"37\u00d750=\"),\n    @(\"75\u00d754=\", \"90\u00d764=\"),\n    @(\"28\u00d755=\", \"17\u00d711=\"),\n    @(\"13\u00d788=\", \"95\u00d724=\"),\n    @(\"16\u00d780=\", \"32\u00d727=\"),\n    @(\"79\u00d791=\", \"19\u00d742=\"),\n    @(\"61\u00d737=\", \"26\u00d772=\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $result = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $result) {\n        throw \"Replacement failed for: $findText\"\n    }\n}\n"}
